$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "English" "Inglés"
Replace-Text " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portugués / Francés / Tailandés / Vietnamita / Español"
Replace-Text "Brief" "Breve"
Replace-Text "An email sent to confirmed attendees of the event. We want to share the flight and accommodation booking details with them." "Un correo electrónico enviado a los asistentes confirmados del evento. Queremos compartir los detalles de la reserva de vuelo y alojamiento con ellos."
Replace-Text "Target audience" "Público objetivo"
Replace-Text "Event attendees" "Asistentes al evento"
Replace-Text "Subject: " "Asunto: "
Replace-Text "Here are your booking details for" "Aquí tiene los datos de su reserva para"
Replace-Text "We can’t wait to meet you! " "¡Estamos impacientes por conocerte! "
Replace-Text "Hi " "Hola "
Replace-Text "We hope you’re as excited as we are for " "Esperamos que estés tan emocionado como nosotros por el "
Replace-Text ". As we’re nearing the event, we’ve made all the preparations to have you with us for this " ". A medida que nos acercamos al evento, hemos hecho todos los preparativos para recibirte en esta "
Replace-Text "conference/seminar/trip" "conferencia/seminario/viaje"
Replace-Text "In this email, we’ve linked/attached the following documents:" "En este correo electrónico, hemos enlazado/adjuntado los siguientes documentos:"
Replace-Text "Your return flight tickets" "Tus billetes de avión de ida y vuelta"
Replace-Text "Your accommodation booking details" "Los datos de tu reserva de alojamiento"
Replace-Text "Your visa information " "La información de su visa "
Replace-Text "(if applicable)" "(si corresponde)"
Replace-Text "If you have any questions, please contact us via " "Si tienes alguna pregunta, entra en contacto con nosotros por "
Replace-Text "If you have any questions, please contact your country manager, " "Si tienes alguna pregunta, entra en contacto con el gestor de tu país "
Replace-Text ", at " ", en "
Replace-Text "See you on the " "¡Nos vemos el "
Replace-Text "[DD]th" "día [DD]"
Replace-Text " or " " o "
Replace-Text " or " " o "

Replace-Text "choose either one" "elija uno de los dos"
Replace-Text "check if these are the documents included" "verifique si estos son los documentos incluidos"
Replace-Text "choose one" "elija uno"
